$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.930305123329163
$ws.Range("B1").Value = 2.753225803375244
$ws.Range("C1").Value = 3.018086910247803
$ws.Range("D1").Value = 2.679993629455566
$ws.Range("E1").Value = 1.002297282218933
